$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 4334
$ws.Range("J2").Value = 5001
$ws.Range("L2").Value = 5001
$ws.Range("N2").Value = -5227
# Row 7
$ws.Range("H7").Value = 13991.8
$ws.Range("I7").Value = 3321.3333
$ws.Range("J7").Value = 29997.5
$ws.Range("K7").Value = 3321.3333
$ws.Range("L7").Value = 29997.5
$ws.Range("M7").Value = -3209.3333
$ws.Range("N7").Value = -30221.5
# Row 14
$ws.Range("H14").Value = 13991.8
$ws.Range("I14").Value = 3321.3333
$ws.Range("J14").Value = 29997.5
$ws.Range("K14").Value = 3321.3333
$ws.Range("L14").Value = 29997.5
$ws.Range("M14").Value = -3130.3333
$ws.Range("N14").Value = -30379.5
# Row 39
$ws.Range("H39").Value = 253.33333
$ws.Range("I39").Value = 253.33333
$ws.Range("K39").Value = 759.99999
$ws.Range("M39").Value = -463.99999
# Row 40
$ws.Range("H40").Value = 1657.2727
$ws.Range("I40").Value = 725
$ws.Range("K40").Value = 725
$ws.Range("M40").Value = -550
# Row 43
$ws.Range("H43").Value = 1111
$ws.Range("I43").Value = 1111
$ws.Range("K43").Value = 1111
$ws.Range("M43").Value = -1042
# Row 107
$ws.Range("H107").Value = 202.375
$ws.Range("I107").Value = 202.85715
$ws.Range("J107").Value = 199
$ws.Range("K107").Value = 202.85715
$ws.Range("L107").Value = 199
$ws.Range("M107").Value = 1717.14285
$ws.Range("N107").Value = -4039
# Row 112
$ws.Range("H112").Value = 1160.725
$ws.Range("J112").Value = 1467.7407
$ws.Range("L112").Value = 4403.2221
$ws.Range("N112").Value = -6619.2221
# Row 137
$ws.Range("H137").Value = 5217.778
$ws.Range("I137").Value = 1989
$ws.Range("J137").Value = 5621.375
$ws.Range("K137").Value = 5967
$ws.Range("L137").Value = 16864.125
$ws.Range("M137").Value = -3417
$ws.Range("N137").Value = -21964.125
# Row 138
$ws.Range("H138").Value = 2715.3276
$ws.Range("I138").Value = 1811.6923
$ws.Range("J138").Value = 2976.3777
$ws.Range("K138").Value = 5435.0769
$ws.Range("L138").Value = 8929.133099999999
$ws.Range("M138").Value = -295.0769
$ws.Range("N138").Value = -19209.1331

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6088
$ws.Range("I32").Value = 6088
$ws.Range("K32").Value = 6088
$ws.Range("M32").Value = -5801
# Row 61
$ws.Range("H61").Value = 2193.4443
$ws.Range("I61").Value = 1609.4
$ws.Range("K61").Value = 1609.4
$ws.Range("M61").Value = -1397.4
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
# Row 110
$ws.Range("H110").Value = 952.3333
$ws.Range("I110").Value = 952.3333
$ws.Range("K110").Value = 952.3333
$ws.Range("M110").Value = 1092.6667
# Row 136
$ws.Range("H136").Value = 2193.4443
$ws.Range("I136").Value = 1609.4
$ws.Range("K136").Value = 4828.200000000001
$ws.Range("M136").Value = -2278.200000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 617.2
$ws.Range("I22").Value = 608.2222
$ws.Range("K22").Value = 608.2222
$ws.Range("M22").Value = -435.2222
# Row 94
$ws.Range("H94").Value = 758.2857
$ws.Range("I94").Value = 384.66666
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 384.66666
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = 66.33334000000002
$ws.Range("N94").Value = -3902
# Row 107
$ws.Range("H107").Value = 1602.8
$ws.Range("I107").Value = 1476.5
$ws.Range("K107").Value = 1476.5
$ws.Range("M107").Value = 443.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2566.5557
$ws.Range("I31").Value = 2726.4285
$ws.Range("J31").Value = 2007
$ws.Range("K31").Value = 2726.4285
$ws.Range("L31").Value = 2007
$ws.Range("M31").Value = -2431.4285
$ws.Range("N31").Value = -2597
# Row 34
$ws.Range("H34").Value = 2566.5557
$ws.Range("I34").Value = 2726.4285
$ws.Range("J34").Value = 2007
$ws.Range("K34").Value = 2726.4285
$ws.Range("L34").Value = 2007
$ws.Range("M34").Value = -2524.4285
$ws.Range("N34").Value = -2411
# Row 107
$ws.Range("H107").Value = 1142.1177
$ws.Range("I107").Value = 523.1667
$ws.Range("J107").Value = 2627.6
$ws.Range("K107").Value = 523.1667
$ws.Range("L107").Value = 2627.6
$ws.Range("M107").Value = 1396.8333
$ws.Range("N107").Value = -6467.6
# Row 132
$ws.Range("H132").Value = 3499.75
$ws.Range("J132").Value = 3571.1428
$ws.Range("L132").Value = 10713.4284
$ws.Range("N132").Value = -15773.4284

$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 18442.455
$ws.Range("I56").Value = 18442.455
$ws.Range("K56").Value = 18442.455
$ws.Range("M56").Value = -17912.455
# Row 107
$ws.Range("H107").Value = 500
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 1500
$ws.Range("N107").Value = -5340
# Row 108
$ws.Range("H108").Value = 2518
$ws.Range("I108").Value = 2518
$ws.Range("K108").Value = 7554
$ws.Range("M108").Value = -4674
# Row 109
$ws.Range("H109").Value = 4299.2354
$ws.Range("I109").Value = 1047.5
$ws.Range("J109").Value = 4732.8
$ws.Range("K109").Value = 3142.5
$ws.Range("L109").Value = 14198.4
$ws.Range("M109").Value = -2102.5
$ws.Range("N109").Value = -16278.4

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3629.1
$ws.Range("I80").Value = 2926.4285
$ws.Range("K80").Value = 2926.4285
$ws.Range("M80").Value = -1928.4285
# Row 83
$ws.Range("H83").Value = 3629.1
$ws.Range("I83").Value = 2926.4285
$ws.Range("K83").Value = 14632.1425
$ws.Range("M83").Value = -9640.1425
# Row 107
$ws.Range("H107").Value = 856.7143
$ws.Range("I107").Value = 856.7143
$ws.Range("K107").Value = 856.7143
$ws.Range("M107").Value = 1063.2857

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 999
$ws.Range("K22").Value = 999
$ws.Range("M22").Value = -704
# Row 27
$ws.Range("H27").Value = 999
$ws.Range("I27").Value = 999
$ws.Range("K27").Value = 999
$ws.Range("M27").Value = -892

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 775.75
$ws.Range("I107").Value = 666.6667
$ws.Range("K107").Value = 2000.0001
$ws.Range("M107").Value = -80.00009999999997
# Row 113
$ws.Range("H113").Value = 886.3125
$ws.Range("I113").Value = 947.1539
$ws.Range("J113").Value = 622.6667
$ws.Range("K113").Value = 2841.4617
$ws.Range("L113").Value = 1868.0001
$ws.Range("M113").Value = -671.4616999999998
$ws.Range("N113").Value = -6208.0001
# Row 132
$ws.Range("H132").Value = 3734.4211
$ws.Range("J132").Value = 3979.5833
$ws.Range("L132").Value = 11938.7499
$ws.Range("N132").Value = -16998.7499
